$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Wed Apr 17 09:30:35 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with the latest scrape.
# Rows 31/32 also swapped rank order (ImmutableX now above Mantle), so
# their Coin/Link/Price/Volume cells are fully replaced.

# Price values that look like plain numbers ("139.40", "78.80", ...) must be
# pinned to Text format first, otherwise Excel would coerce them to numeric
# values and silently drop significant trailing zeros.
$priceTextCells = @(
    "D5",
    "D6",
    "D9",
    "D14",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D37",
    "D39",
    "D40",
    "D42",
    "D43",
    "D44",
    "D46",
    "D47"
)
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.534.83'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.085.80'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("D5").Value = '546.19'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = '139.40'
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.078.30'
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("E10").Value = '  +0.77%  '
$ws.Range("E11").Value = '  +1.66%  '
$ws.Range("E12").Value = '  -3.00%  '
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("D14").Value = '35.11'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").Value = '3.583.65'
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Value = '63.517.43'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").Value = '3.079.81'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("D20").Value = '477.02'
$ws.Range("E20").Value = '  -2.67%  '
$ws.Range("D21").Value = '13.54'
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").Value = '0.703'
$ws.Range("E22").Value = '  -2.45%  '
$ws.Range("D23").Value = '7.10'
$ws.Range("E23").Value = '  -2.78%  '
$ws.Range("D24").Value = '78.80'
$ws.Range("D25").Value = '12.25'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '2.73'
$ws.Range("E27").Value = '  -1.48%  '
$ws.Range("D28").Value = '7.97'
$ws.Range("E28").Value = '  -6.98%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '26.33'
$ws.Range("E30").Value = '  -1.38%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.90'
$ws.Range("E31").Value = '  -3.94%  '
$ws.Range("B32").Value = 'Mantle'
$ws.Range("C32").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D32").Value = '1.16'
$ws.Range("E32").Value = '  +2.67%  '
$ws.Range("D33").Value = '59.12'
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("E34").Value = '  -7.56%  '
$ws.Range("D35").Value = '5.53'
$ws.Range("E35").Value = '  +7.54%  '
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("D37").Value = '490.20'
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("D38").Value = '3.274.11'
$ws.Range("E38").Value = '  +3.76%  '
$ws.Range("D39").Value = '0.0404'
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").Value = '0.0800'
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("D42").Value = '8.19'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = '2.62'
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("D44").Value = '0.254'
$ws.Range("E44").Value = '  -1.81%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = '25.50'
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").Value = '124.34'
$ws.Range("E47").Value = '  +3.02%  '
$ws.Range("E48").Value = '  -2.07%  '
$ws.Range("D49").Value = '0.0₃0532'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("E51").Value = '  -0.37%  '
